$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 112390426
$ws.Range("J24").Value = "fruktkroppar"
$ws.Range("Q24").Value = 509076
$ws.Range("R24").Value = 6783959
$ws.Range("A25").Value = 112390567
$ws.Range("B25").Value = 90806
$ws.Range("E25").Value = 4361
$ws.Range("F25").Value = "Orange taggsvamp"
$ws.Range("G25").Value = "Hydnellum aurantiacum"
$ws.Range("H25").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I25").Value = "19"
$ws.Range("J25").Value = "fruktkroppar"
$ws.Range("Q25").Value = 509010
$ws.Range("R25").Value = 6783836
$ws.Range("A26").Value = 112390451
$ws.Range("B26").Value = 90814
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 4364
$ws.Range("F26").Value = "Dropptaggsvamp"
$ws.Range("G26").Value = "Hydnellum ferrugineum"
$ws.Range("H26").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I26").NumberFormat = "@"
$ws.Range("I26").Value = "3"
$ws.Range("Q26").Value = 509076
$ws.Range("R26").Value = 6783959
$ws.Range("AH26").Value = "Sandtallskog"
$ws.Range("A27").Value = 112390119
$ws.Range("B27").Value = 90830
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 2059
$ws.Range("F27").Value = "Skrovlig taggsvamp"
$ws.Range("G27").Value = "Hydnellum scabrosum"
$ws.Range("H27").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q27").Value = 509093
$ws.Range("R27").Value = 6784215
$ws.Range("A28").Value = 112390524
$ws.Range("B28").Value = 90814
$ws.Range("E28").Value = 4364
$ws.Range("F28").Value = "Dropptaggsvamp"
$ws.Range("G28").Value = "Hydnellum ferrugineum"
$ws.Range("H28").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I28").Value = ""
$ws.Range("J28").Value = ""
$ws.Range("Q28").Value = 509060
$ws.Range("R28").Value = 6783866
$ws.Range("AH28").Value = ""
$ws.Range("A29").Value = 112390630
$ws.Range("B29").Value = 90857
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 5448
$ws.Range("F29").Value = "Svartvit taggsvamp"
$ws.Range("G29").Value = "Phellodon connatus"
$ws.Range("H29").Value = "(Schultz) nom.prov"
$ws.Range("Q29").Value = 509014
$ws.Range("R29").Value = 6783848
$ws.Range("A30").Value = 112390287
$ws.Range("I30").Value = ""
$ws.Range("J30").Value = ""
$ws.Range("Q30").Value = 509070
$ws.Range("R30").Value = 6784097
$ws.Range("A31").Value = 112390262
$ws.Range("B31").Value = 90808
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 4362
$ws.Range("F31").Value = "Blå taggsvamp"
$ws.Range("G31").Value = "Hydnellum caeruleum"
$ws.Range("H31").Value = "(Hornem.) P.Karst."
$ws.Range("I31").Value = ""
$ws.Range("Q31").Value = 509072
$ws.Range("R31").Value = 6784116
$ws.Range("A32").Value = 112389988
$ws.Range("B32").Value = 89072
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 256703
$ws.Range("F32").Value = "Tallfingersvamp"
$ws.Range("G32").Value = "Ramaria eosanguinea"
$ws.Range("H32").Value = "R.H.Petersen"
$ws.Range("J32").Value = ""
$ws.Range("Q32").Value = 509101
$ws.Range("R32").Value = 6784234
$ws.Range("A33").Value = 112390292
$ws.Range("B33").Value = 90808
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 4362
$ws.Range("F33").Value = "Blå taggsvamp"
$ws.Range("G33").Value = "Hydnellum caeruleum"
$ws.Range("H33").Value = "(Hornem.) P.Karst."
$ws.Range("I33").NumberFormat = "@"
$ws.Range("I33").Value = "1"
$ws.Range("Q33").Value = 509065
$ws.Range("R33").Value = 6784066
$ws.Range("AH33").Value = ""
$ws.Range("A34").Value = 112390382
$ws.Range("B34").Value = 90830
$ws.Range("E34").Value = 2059
$ws.Range("F34").Value = "Skrovlig taggsvamp"
$ws.Range("G34").Value = "Hydnellum scabrosum"
$ws.Range("H34").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I34").Value = ""
$ws.Range("J34").Value = ""
$ws.Range("Q34").Value = 509061
$ws.Range("R34").Value = 6784061
$ws.Range("A35").Value = 112390398
$ws.Range("Q35").Value = 509066
$ws.Range("R35").Value = 6784010
$ws.Range("A36").Value = 112390509
$ws.Range("Q36").Value = 509056
$ws.Range("R36").Value = 6783885
$ws.Range("AH36").Value = ""
$ws.Range("A37").Value = 112390256
$ws.Range("B37").Value = 90448
$ws.Range("E37").Value = 4745
$ws.Range("F37").Value = "Tallriska"
$ws.Range("G37").Value = "Lactarius musteus"
$ws.Range("H37").Value = "Fr."
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "1"
$ws.Range("J37").Value = "fruktkroppar"
$ws.Range("Q37").Value = 509090
$ws.Range("R37").Value = 6784191
$ws.Range("AH37").Value = "Sandtallskog"
$ws.Range("A38").Value = 112390031
$ws.Range("B38").Value = 90826
$ws.Range("D38").Value = "LC"
$ws.Range("E38").Value = 4366
$ws.Range("F38").Value = "Skarp dropptaggsvamp"
$ws.Range("G38").Value = "Hydnellum peckii"
$ws.Range("H38").Value = "Banker"
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = "1"
$ws.Range("Q38").Value = 509098
$ws.Range("R38").Value = 6784229
$ws.Range("AH38").Value = "Sandtallskog"
$ws.Range("A39").Value = 112389959
$ws.Range("B39").Value = 89072
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 256703
$ws.Range("F39").Value = "Tallfingersvamp"
$ws.Range("G39").Value = "Ramaria eosanguinea"
$ws.Range("H39").Value = "R.H.Petersen"
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "1"
$ws.Range("J39").Value = "fruktkroppar"
$ws.Range("Q39").Value = 509111
$ws.Range("R39").Value = 6784257
